$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell values (shared strings change: name2L/nameS /x -> test group/nameB/z)
$ws.Range("A1").Value = "test group"
$ws.Range("A2").Value = "nameB"
$ws.Range("A3").Value = "z"

# Widen columns A:C to a custom width of 26 characters
$ws.Columns("A:C").ColumnWidth = 25.17

# Move the active selection to B2
$ws.Range("B2").Select()
